$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2698.8667
$ws.Range("J70").Value = 2607.1667
$ws.Range("L70").Value = 7821.500100000001
$ws.Range("N70").Value = -8361.500100000001
$ws.Range("H73").Value = 2698.8667
$ws.Range("J73").Value = 2607.1667
$ws.Range("L73").Value = 7821.500100000001
$ws.Range("N73").Value = -9693.500100000001
$ws.Range("H86").Value = 2151.1365
$ws.Range("I86").Value = 2318.111
$ws.Range("K86").Value = 2318.111
$ws.Range("M86").Value = -1195.111
$ws.Range("H89").Value = 2151.1365
$ws.Range("I89").Value = 2318.111
$ws.Range("K89").Value = 11590.555
$ws.Range("M89").Value = -5974.555
$ws.Range("H134").Value = 96759.664
$ws.Range("J134").Value = 96759.664
$ws.Range("L134").Value = 96759.664
$ws.Range("N134").Value = -106899.664
$ws.Range("H135").Value = 481.04544
$ws.Range("I135").Value = 444.83334
$ws.Range("K135").Value = 4003.50006
$ws.Range("M135").Value = -1468.50006
$ws.Range("H138").Value = 2824.7144
$ws.Range("I138").Value = 1932
$ws.Range("J138").Value = 3181.8
$ws.Range("K138").Value = 5796
$ws.Range("L138").Value = 9545.400000000001
$ws.Range("M138").Value = -656
$ws.Range("N138").Value = -19825.4
$ws.Range("H141").Value = 1676.909
$ws.Range("I141").Value = 1707.742
$ws.Range("K141").Value = 5123.226
$ws.Range("M141").Value = 56.77400000000034

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3044.6868
$ws.Range("I32").Value = 3044.6868
$ws.Range("K32").Value = 3044.6868
$ws.Range("M32").Value = -2757.6868
$ws.Range("H61").Value = 14180.923
$ws.Range("I61").Value = 9037.105
$ws.Range("K61").Value = 9037.105
$ws.Range("M61").Value = -8825.105
$ws.Range("H102").Value = 4027.3572
$ws.Range("I102").Value = 2716.6365
$ws.Range("J102").Value = 8833.333000000001
$ws.Range("K102").Value = 2716.6365
$ws.Range("L102").Value = 8833.333000000001
$ws.Range("M102").Value = -1094.6365
$ws.Range("N102").Value = -12077.333
$ws.Range("H132").Value = 4416.05
$ws.Range("I132").Value = 4385.316
$ws.Range("K132").Value = 13155.948
$ws.Range("M132").Value = -10625.948
$ws.Range("H136").Value = 14180.923
$ws.Range("I136").Value = 9037.105
$ws.Range("K136").Value = 27111.315
$ws.Range("M136").Value = -24561.315

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 62813
$ws.Range("J2").Value = 69219.5
$ws.Range("L2").Value = 69219.5
$ws.Range("N2").Value = -69445.5
$ws.Range("H20").Value = 3173.6428
$ws.Range("I20").Value = 2957.158
$ws.Range("J20").Value = 3630.6667
$ws.Range("K20").Value = 2957.158
$ws.Range("L20").Value = 3630.6667
$ws.Range("M20").Value = -2710.158
$ws.Range("N20").Value = -4124.6667
$ws.Range("H105").Value = 1877
$ws.Range("I105").Value = 1995
$ws.Range("J105").Value = 1798.3334
$ws.Range("K105").Value = 1995
$ws.Range("L105").Value = 1798.3334
$ws.Range("M105").Value = -248
$ws.Range("N105").Value = -5292.3334
$ws.Range("H134").Value = 2640.1929
$ws.Range("I134").Value = 2614.309
$ws.Range("K134").Value = 7842.927000000001
$ws.Range("M134").Value = -5307.927000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7409.857
$ws.Range("I31").Value = 5399.6
$ws.Range("J31").Value = 8526.666999999999
$ws.Range("K31").Value = 5399.6
$ws.Range("L31").Value = 8526.666999999999
$ws.Range("M31").Value = -5104.6
$ws.Range("N31").Value = -9116.666999999999
$ws.Range("H34").Value = 7409.857
$ws.Range("I34").Value = 5399.6
$ws.Range("J34").Value = 8526.666999999999
$ws.Range("K34").Value = 5399.6
$ws.Range("L34").Value = 8526.666999999999
$ws.Range("M34").Value = -5197.6
$ws.Range("N34").Value = -8930.666999999999
$ws.Range("H58").Value = 6479.815
$ws.Range("I58").Value = 4691
$ws.Range("K58").Value = 4691
$ws.Range("M58").Value = -4488
$ws.Range("H74").Value = 42717
$ws.Range("J74").Value = 42717
$ws.Range("L74").Value = 42717
$ws.Range("N74").Value = -44465
$ws.Range("H77").Value = 42717
$ws.Range("J77").Value = 42717
$ws.Range("L77").Value = 128151
$ws.Range("N77").Value = -136887
$ws.Range("H86").Value = 46438.77
$ws.Range("I86").Value = 6910.8
$ws.Range("K86").Value = 6910.8
$ws.Range("M86").Value = -5787.8
$ws.Range("H89").Value = 46438.77
$ws.Range("I89").Value = 6910.8
$ws.Range("K89").Value = 34554
$ws.Range("M89").Value = -28938
$ws.Range("H103").Value = 74288
$ws.Range("I103").Value = 70000
$ws.Range("K103").Value = 70000
$ws.Range("M103").Value = -68828
$ws.Range("H132").Value = 4518.1
$ws.Range("I132").Value = 4680.1177
$ws.Range("K132").Value = 14040.3531
$ws.Range("M132").Value = -11510.3531
$ws.Range("H136").Value = 6479.815
$ws.Range("I136").Value = 4691
$ws.Range("K136").Value = 14073
$ws.Range("M136").Value = -11523

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 15349.5
$ws.Range("I123").Value = 4933
$ws.Range("J123").Value = 21599.4
$ws.Range("K123").Value = 14799
$ws.Range("L123").Value = 64798.2
$ws.Range("M123").Value = -12349
$ws.Range("N123").Value = -69698.20000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14582.934
$ws.Range("I70").Value = 10963.714
$ws.Range("J70").Value = 17749.75
$ws.Range("K70").Value = 10963.714
$ws.Range("L70").Value = 17749.75
$ws.Range("M70").Value = -10693.714
$ws.Range("N70").Value = -18289.75
$ws.Range("H73").Value = 14582.934
$ws.Range("I73").Value = 10963.714
$ws.Range("J73").Value = 17749.75
$ws.Range("K73").Value = 10963.714
$ws.Range("L73").Value = 17749.75
$ws.Range("M73").Value = -10027.714
$ws.Range("N73").Value = -19621.75
$ws.Range("H96").Value = 49425.6
$ws.Range("J96").Value = 49425.6
$ws.Range("L96").Value = 49425.6
$ws.Range("N96").Value = -54917.6
$ws.Range("H99").Value = 8278.888999999999
$ws.Range("I99").Value = 7438.75
$ws.Range("J99").Value = 15000
$ws.Range("K99").Value = 7438.75
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -5192.75
$ws.Range("N99").Value = -19492
$ws.Range("H102").Value = 4092.087
$ws.Range("I102").Value = 1882.5625
$ws.Range("K102").Value = 1882.5625
$ws.Range("M102").Value = -260.5625
$ws.Range("H132").Value = 1821.7142
$ws.Range("I132").Value = 1612.8
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 4838.4
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -2308.4
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H82").Value = 1222.1818
$ws.Range("I82").Value = 834.8570999999999
$ws.Range("J82").Value = 1900
$ws.Range("K82").Value = 834.8570999999999
$ws.Range("L82").Value = 1900
$ws.Range("M82").Value = -473.8570999999999
$ws.Range("N82").Value = -2622
$ws.Range("H85").Value = 1222.1818
$ws.Range("I85").Value = 834.8570999999999
$ws.Range("J85").Value = 1900
$ws.Range("K85").Value = 834.8570999999999
$ws.Range("L85").Value = 1900
$ws.Range("M85").Value = 413.1429000000001
$ws.Range("N85").Value = -4396
$ws.Range("H93").Value = 6610.1665
$ws.Range("I93").Value = 1728.5172
$ws.Range("K93").Value = 1728.5172
$ws.Range("M93").Value = -480.5172
$ws.Range("H101").Value = 88333.336
$ws.Range("J101").Value = 88333.336
$ws.Range("L101").Value = 88333.336
$ws.Range("N101").Value = -94823.336
$ws.Range("H119").Value = 54200
$ws.Range("J119").Value = 54200
$ws.Range("L119").Value = 54200
$ws.Range("N119").Value = -63876
$ws.Range("H132").Value = 8584.305
$ws.Range("I132").Value = 8808.319
$ws.Range("K132").Value = 26424.957
$ws.Range("M132").Value = -23894.957
$ws.Range("H136").Value = 3650.762
$ws.Range("I136").Value = 3469.3948
$ws.Range("J136").Value = 5373.75
$ws.Range("K136").Value = 10408.1844
$ws.Range("L136").Value = 16121.25
$ws.Range("M136").Value = -7858.1844
$ws.Range("N136").Value = -21221.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 149985
$ws.Range("J5").Value = 149985
$ws.Range("L5").Value = 149985
$ws.Range("N5").Value = -150209
$ws.Range("H103").Value = 47250
$ws.Range("J103").Value = 47250
$ws.Range("L103").Value = 47250
$ws.Range("N103").Value = -49594
$ws.Range("H132").Value = 2379.8438
$ws.Range("I132").Value = 2379.8438
$ws.Range("K132").Value = 7139.5314
$ws.Range("M132").Value = -4609.5314
$ws.Range("H136").Value = 3666.2415
$ws.Range("I136").Value = 2931.75
$ws.Range("K136").Value = 8795.25
$ws.Range("M136").Value = -6245.25
